# Add season-record columns (Wins / Losses / Ties) to the KCR 2002 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same bold/bordered/centered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-54): 62 wins, 100 losses, 0 ties
$ws.Range("AD2:AD54").Value = 62
$ws.Range("AE2:AE54").Value = 100
$ws.Range("AF2:AF54").Value = 0
